$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell selection to E7 (was F7)
$ws.Range("E7").Select()

# Row 7: apply the "thick bottom border" row styling (matches rows 3)
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Range("A7:B7").Borders.Item(9).LineStyle = 1   # xlEdgeBottom (9) - double line
$ws.Range("A7:B7").Borders.Item(9).Weight = -4138  # xlThick
$ws.Range("A7:B7").Borders.Item(9).ColorIndex = 1

# Row 8: apply the "thick top + bottom border" row styling (matches row 9)
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Range("A8:B8").Borders.Item(8).LineStyle = 1   # xlEdgeTop (8)
$ws.Range("A8:B8").Borders.Item(8).Weight = -4138
$ws.Range("A8:B8").Borders.Item(8).ColorIndex = 1
$ws.Range("A8:B8").Borders.Item(9).LineStyle = 1   # xlEdgeBottom (9)
$ws.Range("A8:B8").Borders.Item(9).Weight = -4138
$ws.Range("A8:B8").Borders.Item(9).ColorIndex = 1

# B8: change status from Pending to In Progress, and apply the "In Progress" cell style
$ws.Range("B8").Value = "In Progress"
$ws.Range("B8").Style = "Check Cell"
